$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hardcoded CAR_SERIES range-check formulas to use the new
# exclude() based function syntax.
$ws.Range("E4").Value = 'exclude(CAR_SERIES ,"MILW") <= 120209'
$ws.Range("D4").Value = 'exclude(CAR_SERIES ,"MILW" )>= 120000 '

# Widen column D slightly to fit the new, longer formula text
# (closest attainable width to the target 52.140625 given this engine's
# internal pixel-snapping of column widths).
$ws.Columns.Item(4).ColumnWidth = 51.333333

# Move the active selection to E4 to match the author's final cursor position.
$ws.Range("E4").Select() | Out-Null
